$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in row 6 with the new resource (REC160: "La fecundación in vitro"),
# matching the pattern used by rows 2-5 (Area / Codigo guion / Numero recurso /
# Nombre recurso / Matricula video / Tipo / Origen / Observaciones).
$ws.Range("A6").Value = "Ciencias Naturales"
$ws.Range("B6").Value = "CN_08_05_CO"
$ws.Range("C6").Value = 16
$ws.Range("D6").Value = "La fecundación in vitro"
$ws.Range("F6").Value = "Recurso"
$ws.Range("G6").Value = "AP"
$ws.Range("H6").Value = "Cambiar acento"

# The row now wraps onto two lines like the rows above it.
$ws.Range("A6:H6").RowHeight = 30

# Move the active selection to E2 (matricula video column on the first data row).
[void]$ws.Range("E2").Select()
